$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cryptos list refresh: Price (D) / Volume(1h) (E) columns for rows 2-51.
# Source cells are literal text (e.g. "67.580.78", "  -1.83%  "). A couple
# of new Price values ("551.90", "159.60") are valid decimals with a
# trailing zero, which Excel's value-assignment auto-detects as numbers and
# would silently round to "551.9" / "159.6" — prefix those with a leading
# apostrophe so they stay text, matching the source file's stored values.

$ws.Range("D2").Value = "67.580.78"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").Value = "2.426.97"
$ws.Range("E3").Value = "  -1.47%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'551.90"
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("D6").Value = "'159.60"
$ws.Range("E6").Value = "  -1.48%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.507"
$ws.Range("E8").Value = "  +0.69%  "
$ws.Range("D9").Value = "0.158"
$ws.Range("E9").Value = "  +5.81%  "
$ws.Range("E10").Value = "  -0.91%  "
$ws.Range("D11").Value = "0.328"
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("D13").Value = "67.517.97"
$ws.Range("E13").Value = "  -1.52%  "
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "22.95"
$ws.Range("E15").Value = "  -2.26%  "
$ws.Range("D16").Value = "10.34"
$ws.Range("E16").Value = "  -3.61%  "
$ws.Range("D17").Value = "328.87"
$ws.Range("E17").Value = "  -3.72%  "
$ws.Range("D18").Value = "6.84"
$ws.Range("E18").Value = "  -2.72%  "
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "1.84"
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("D22").Value = "65.96"
$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("E23").Value = "  -1.65%  "
$ws.Range("D24").Value = "8.06"
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("D25").Value = "0.0₃0803"
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("D26").Value = "7.01"
$ws.Range("E26").Value = "  -2.12%  "
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "415.76"
$ws.Range("E28").Value = "  -4.76%  "
$ws.Range("D29").Value = "1.12"
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("E30").Value = "  -1.25%  "
$ws.Range("D31").Value = "160.21"
$ws.Range("E31").Value = "  +1.95%  "
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").Value = "17.77"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("E35").Value = "  -3.76%  "
$ws.Range("D36").Value = "0.294"
$ws.Range("E36").Value = "  -2.56%  "
$ws.Range("E37").Value = "  -4.18%  "
$ws.Range("D38").Value = "1.45"
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("E39").Value = "  -3.01%  "
$ws.Range("E40").Value = "  -3.35%  "
$ws.Range("D41").Value = "3.31"
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("D42").Value = "129.99"
$ws.Range("E42").Value = "  -1.96%  "
$ws.Range("D43").Value = "0.0706"
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("D45").Value = "0.553"
$ws.Range("E45").Value = "  -1.54%  "
$ws.Range("D46").Value = "0.0913"
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("E48").Value = "  -7.52%  "
$ws.Range("D49").Value = "16.47"
$ws.Range("E49").Value = "  -1.94%  "
$ws.Range("D50").Value = "0.0₆0203"
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("E51").Value = "  -0.25%  "
